# ADD results from server
#
# The workbook has 6 sheets (one per model year: 2025, 2030, 2035, 2040,
# 2045, 2050), each with the same 15 technology-code column headers in
# row 1 and a single data row (row 2) below them.
#
# The server re-ran the model and returned an updated column layout
# (two new technologies - "gb" and "btes" - are now reported, while "gt"
# and "dgt" are no longer present) together with refreshed investment
# cost figures for every sheet.

$wb = $excel.ActiveWorkbook

$headers = @("eb","gb","hp","st","wi","ieh","chp","ac","ab_ct","ab_hp","cp_ct","cp_hp","ttes","btes","ites")

# New row-2 values (in column order matching $headers above) for each of
# the six sheets, keyed by sheet index (1-based, matching tab order
# 2025..2050).
$sheetValues = @{
    1 = @(39063.99109145206, 0, 483537.6274462014, 0, 2897240.114301849, 94331.34471502228, 0, 25342.77928792104, 0, 0, 0, 0, 0, 23638.06126801545, 19940.13531829346)
    2 = @(30846.52922536713, 0, 1495599.874611417, 0, 0, 70193.79982138964, 0, 56602.42752520426, 0, 0, 0, 0, 0, 51649.16401227913, 42574.77934331147)
    3 = @(242452.4252219552, 0, 943335.270081223, 0, 0, 1425.925979620855, 0, 39373.98526588717, 0, 0, 0, 0, 0, 53308.16490721726, 30023.09380555204)
    4 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 11578.49752443177, 0)
    5 = @(76705.58894163162, 1930.947398408091, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 28147.3462746636, 8312.661449003012)
    6 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Row 1: refreshed header / column codes
    for ($i = 0; $i -lt $headers.Count; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }

    # Row 2: refreshed investment-cost data for this sheet
    $values = $sheetValues[$s]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item(2, $i + 1).Value = $values[$i]
    }
}
